# Update the cryptos worksheet with the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain TEXT (no numeric/date auto-coercion), and
# strip the temporary "Text" number-format back off afterwards so the
# cell's style stays exactly as it started (no explicit style index).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# --- Rows whose Price (D) and/or Volume(1h) (E) values changed, but Coin/Link stayed the same ---

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "41.774.88"
$ws.Range("E2").Value = "  +2.24%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "2.263.63"
$ws.Range("E3").Value = "  +0.99%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.02%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "303.94"
$ws.Range("E5").Value = "  -0.02%  "

# Row 6 - Solana
Set-TextValue $ws.Range("D6") "91.99"
$ws.Range("E6").Value = "  +0.67%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  +1.90%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.12%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  -0.07%  "

# Row 10 - Avalanche
$ws.Range("E10").Value = "  +0.97%  "

# Row 11 - OKB
$ws.Range("E11").Value = "  +1.36%  "

# Row 12 - Dogecoin
$ws.Range("E12").Value = "  +0.95%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  +0.01%  "

# Row 14 - Polkadot
Set-TextValue $ws.Range("D14") "6.61"
$ws.Range("E14").Value = "  +0.90%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D15") "2.614.14"
$ws.Range("E15").Value = "  +1.04%  "

# Row 17 - WrappedEther
Set-TextValue $ws.Range("D17") "2.269.06"
$ws.Range("E17").Value = "  +1.69%  "

# Row 18 - Polygon
Set-TextValue $ws.Range("D18") "0.764"
$ws.Range("E18").Value = "  +1.95%  "

# Row 19 - WrappedBTC
Set-TextValue $ws.Range("D19") "41.667.60"
$ws.Range("E19").Value = "  +2.24%  "

# Row 20 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D20") "12.51"
$ws.Range("E20").Value = "  +6.80%  "

# Row 21 - ShibaInu
Set-TextValue $ws.Range("D21") "0.0$([char]8323)0903"
$ws.Range("E21").Value = "  +0.42%  "

# Row 23 - Litecoin
Set-TextValue $ws.Range("D23") "66.85"
$ws.Range("E23").Value = "  +0.93%  "

# Row 24 - BitcoinCash
Set-TextValue $ws.Range("D24") "239.62"
$ws.Range("E24").Value = "  +0.26%  "

# Row 25 - PancakeSwap
$ws.Range("E25").Value = "  +1.34%  "

# Row 27 - ImmutableX
$ws.Range("E27").Value = "  +3.19%  "

# Row 28 - EthereumClassic
Set-TextValue $ws.Range("D28") "23.99"
$ws.Range("E28").Value = "  +0.27%  "

# Row 29 - Cosmos
$ws.Range("E29").Value = "  +0.50%  "

# Row 30 - Toncoin
$ws.Range("E30").Value = "  -4.54%  "

# Row 31 - Monero
Set-TextValue $ws.Range("D31") "160.55"
$ws.Range("E31").Value = "  +1.58%  "

# Row 32 - InjectiveProtocol
Set-TextValue $ws.Range("D32") "34.39"
$ws.Range("E32").Value = "  +3.62%  "

# Row 33 - Filecoin
$ws.Range("E33").Value = "  +3.89%  "

# Row 34 - FirstDigitalUSD
$ws.Range("E34").Value = "  -0.15%  "

# Row 35 - Hedera
Set-TextValue $ws.Range("D35") "0.0744"
$ws.Range("E35").Value = "  +2.41%  "

# --- Rows 37-40 and 49-51: coins re-ranked (rows swapped/rotated), with updated prices/volumes ---

# Row 37 - now Celestia (was WEMIXToken)
$ws.Range("B37").Value = "Celestia"
$ws.Range("C37").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextValue $ws.Range("D37") "16.89"
$ws.Range("E37").Value = "  +3.85%  "

# Row 38 - now WEMIXToken (was Celestia)
$ws.Range("B38").Value = "WEMIXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Range("D38") "2.38"
$ws.Range("E38").Value = "  +1.36%  "

# Row 39 - now Kaspa (was Stellar)
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws.Range("D39") "0.105"
$ws.Range("E39").Value = "  +0.41%  "

# Row 40 - now Stellar (was Kaspa)
$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws.Range("D40") "0.116"
$ws.Range("E40").Value = "  +1.41%  "

# Row 41 - ARBITRUM
$ws.Range("E41").Value = "  +0.60%  "

# Row 42 - RenderToken
Set-TextValue $ws.Range("D42") "3.95"
$ws.Range("E42").Value = "  +1.44%  "

# Row 43 - Maker
Set-TextValue $ws.Range("D43") "2.025.61"
$ws.Range("E43").Value = "  -3.32%  "

# Row 44 - EnergySwap
Set-TextValue $ws.Range("D44") "19.11"
$ws.Range("E44").Value = "  -3.60%  "

# Row 45 - VeChain
$ws.Range("E45").Value = "  +0.79%  "

# Row 46 - FraxShare
Set-TextValue $ws.Range("D46") "10.35"
$ws.Range("E46").Value = "  +2.19%  "

# Row 47 - ApeXProtocol
$ws.Range("E47").Value = "  +15.13%  "

# Row 48 - NEARProtocol
$ws.Range("E48").Value = "  -1.63%  "

# Row 49 - now Stacks (was TrustWalletToken)
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Range("D49") "1.52"
$ws.Range("E49").Value = "  -0.05%  "

# Row 50 - now BitcoinSV (was Stacks)
$ws.Range("B50").Value = "BitcoinSV"
$ws.Range("C50").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
Set-TextValue $ws.Range("D50") "72.60"
$ws.Range("E50").Value = "  +4.55%  "

# Row 51 - now TrustWalletToken (was BitcoinSV)
$ws.Range("B51").Value = "TrustWalletToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Range("D51") "1.15"
$ws.Range("E51").Value = "  +1.20%  "
